{"js": "// Bold/italic/blue-ify the \"retrain.py\" / \"test.py\" filenames that follow\n// \"Run \" in the two numbered-list steps and their matching bold section\n// headers, then move the hidden \"_GoBack\" bookmark from the blank\n// paragraph after step 6's Google Drive/OneDrive note down to the blank\n// paragraph that now follows the \"7) Run test.py\" bold section header.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/font/bold\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- 1) Split \"...Run retrain.py\" / \"...Run test.py\" runs, formatting just\n//        the filename as bold + italic + blue (#0000FF). Matches both the\n//        plain numbered-list entries (\"6) Run retrain.py\") and the bold,\n//        larger section headers (\"6) Run retrain.py\" @ 12pt bold).\nfor (const p of items) {\n  const m = p.text.match(/^\\d\\)? ?Run (retrain|test)\\.py$/);\n  if (!m) continue;\n  const fileName = `${m[1]}.py`;\n  const results = p.search(fileName, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) continue;\n  const target = results.items[0];\n  target.font.bold = true;\n  target.font.italic = true;\n  target.font.color = \"#0000FF\";\n}\nawait context.sync();\n\n// --- 2) Relocate the \"_GoBack\" bookmark to the blank paragraph right after\n//        the bold \"7) Run test.py\" section header.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nlet afterBoldTestPy = null;\nfor (let i = 0; i < items.length; i++) {\n  if (/^\\d\\)? ?Run test\\.py$/.test(items[i].text) && items[i].font.bold) {\n    afterBoldTestPy = items[i + 1];\n  }\n}\nif (afterBoldTestPy) {\n  afterBoldTestPy.getRange().insertBookmark(\"_GoBack\");\n}\nawait context.sync();\n", "ps1": "# Bold/italic/blue-ify the \"retrain.py\" / \"test.py\" filenames that follow\n# \"Run \" in the two numbered-list steps and their corresponding bold\n# section headers, then move the hidden \"_GoBack\" bookmark from the blank\n# paragraph after step 6's \"Google Drive/OneDrive\" note down to the blank\n# paragraph that now follows the \"7) Run test.py\" section header (this is\n# where Word leaves it after the most recent edit).\n\n$d = $word.ActiveDocument\n\n# --- 1) Split \"...Run retrain.py\" / \"...Run test.py\" runs, formatting the\n#        filename itself as bold+italic+blue (RGB 0,0,255 -> wdColor\n#        16711680 since Word stores colors as BGR). This matches both the\n#        plain numbered-list entries (\"6) Run retrain.py\") and the bold,\n#        larger section headers (\"6) Run retrain.py\" @ 12pt bold).\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -match '^\\d\\)? ?Run (retrain|test)\\.py\\r$') {\n        $target = $matches[1] + \".py\"\n        $pStart = $p.Range.Start\n        $idx = $t.IndexOf($target)\n        if ($idx -ge 0) {\n            $subStart = $pStart + $idx\n            $subEnd = $subStart + $target.Length\n            $rng = $d.Range($subStart, $subEnd)\n            $rng.Font.Bold = 1\n            $rng.Font.Italic = 1\n            $rng.Font.Color = 16711680\n        }\n    }\n}\n\n# --- 2) Relocate the \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$target2 = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -match '^\\d\\)? ?Run test\\.py\\r$' -and $p.Range.Font.Bold) {\n        $target2 = $p\n    }\n}\nif ($target2 -ne $null) {\n    $nextPara = $target2.Next()\n    $d.Bookmarks.Add(\"_GoBack\", $nextPara.Range)\n}\n"}
